# Rename worksheets to reflect rerun model summaries (no urban landuse)
$wb = $excel.ActiveWorkbook

$oldNames = @(
    "summ59709076",
    "summ00390940",
    "summ01058665",
    "summ01704018",
    "summ02302352",
    "summ02950631",
    "summ03561454",
    "summ04178636",
    "summ04846167"
)

$newNames = @(
    "summ05602793",
    "summ06187620",
    "summ06825017",
    "summ07435972",
    "summ08023952",
    "summ08598360",
    "summ09181563",
    "summ09756287",
    "summ10360341"
)

for ($i = 0; $i -lt $oldNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($oldNames[$i])
    $ws.Name = $newNames[$i]
}
